$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D1").Value = "25_02_2024"
$ws.Range("D2").Value = 3414
$ws.Range("D3").Value = 2769
$ws.Range("D4").Value = 3919
$ws.Range("D5").Value = 7058
$ws.Range("D6").Value = 137

$ws.Range("D6").Select()
